# MonteCarloResultsRBMCp214.xlsx edit: added DERS (EENS confidence-interval
# columns) and refreshed the simulated metric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells for the two added columns -----------------------
$ws.Range("Q1").Value = "EENS 95% CI"
$ws.Range("R1").Value = "EENS 99% CI"

# Match the bold / bordered / centered style used by the rest of the
# header row (copy formatting only, values were already set above).
$ws.Range("B1").Copy()
$ws.Range("Q1:R1").PasteSpecial(-4122)

# --- Row 2 (feeder A) refreshed simulation output ----------------------
$ws.Range("F2").Value = 1.630364252423214
$ws.Range("G2").Value = 4428
$ws.Range("H2").Value = 1.203258892709816
$ws.Range("I2").Value = 1.354957160342717
$ws.Range("J2").Value = 284.5410036719707
$ws.Range("K2").Value = 342.3764930088751
$ws.Range("L2").Value = 252.6843674690613
$ws.Range("M2").Value = 0.8722448750464198

# --- Row 3 (feeder B) refreshed simulation output ----------------------
$ws.Range("F3").Value = 1.649225383879152
$ws.Range("G3").Value = 3539
$ws.Range("H3").Value = 1.522935449142997
$ws.Range("I3").Value = 1.082925336597307
$ws.Range("J3").Value = 227.4143206854345
$ws.Range("K3").Value = 346.3373306146219
$ws.Range("L3").Value = 319.8164443200295
$ws.Range("M3").Value = 0.8823355803753462

# --- Row 4 (feeder C) refreshed simulation output ----------------------
$ws.Range("F4").Value = 0.9449432992188636
$ws.Range("G4").Value = 2747
$ws.Range("H4").Value = 1.124162614432925
$ws.Range("I4").Value = 0.8405752753977969
$ws.Range("J4").Value = 176.5208078335373
$ws.Range("K4").Value = 198.4380928359614
$ws.Range("L4").Value = 236.0741490309143
$ws.Range("M4").Value = 0.5055446650820921

# --- Row 5 (TOTAL) refreshed simulation output + new CI columns --------
$ws.Range("J5").Value = 1.092819257445941
$ws.Range("K5").Value = 1.408177645173743
$ws.Range("L5").Value = 1.288573234392699
$ws.Range("M5").Value = 2.260125120503858
$ws.Range("N5").Value = 3268
$ws.Range("P5").Value = 0.01982334740398597
$ws.Range("Q5").Value = "(2.1716844655288594, 2.3472625525421384)"
$ws.Range("R5").Value = "(2.144093623283916, 2.374853394787082)"
